$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph "dorobiť funkciu get_user_id ... - DONE"
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("dorobi")) {
        $targetPara = $p
        break
    }
}

# ------------------------------------------------------------------
# 2. Remove the existing "_GoBack" bookmark (it currently wraps the
#    "DONE" at the end of this paragraph; it will be re-created later
#    inside the newly inserted paragraph).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3. Change the " - " right before the trailing "DONE" into " – "
#    (space, en dash, space) while keeping "DONE" itself bold.
# ------------------------------------------------------------------
$pText = $targetPara.Range.Text
$dashPos = $pText.LastIndexOf(" - ")
$dashStart = $targetPara.Range.Start + $dashPos
$dashEnd = $dashStart + 3
$dashRange = $d.Range($dashStart, $dashEnd)
$dashRange.Text = " " + [char]0x2013 + " "

# ------------------------------------------------------------------
# 4. Insert a brand new list paragraph right after it. It inherits
#    the ListParagraph style / numbering / yellow highlight from the
#    paragraph mark it is split off from.
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("dorobi")) {
        $targetPara = $p
        break
    }
}
$targetPara.Range.InsertParagraphAfter()

$newPara = $null
$afterFound = $false
foreach ($p in $d.Paragraphs) {
    if ($afterFound) {
        $newPara = $p
        break
    }
    if ($p.Range.Text.Contains("dorobi")) {
        $afterFound = $true
    }
}

$newParaStart = $newPara.Range.Start
$newText = "Prerobit priority ciselne (strojove) na pisane (ludske) - DONE"
$newPara.Range.InsertBefore($newText)

# ------------------------------------------------------------------
# 5. Fix up character formatting: the whole line is regular weight
#    except the trailing "DONE", which stays bold.
# ------------------------------------------------------------------
$lineRange = $d.Range($newParaStart, $newParaStart + $newText.Length)
$lineRange.Font.Bold = 0

$doneRange = $d.Range($newParaStart + $newText.Length - 4, $newParaStart + $newText.Length)
$doneRange.Font.Bold = 1

# ------------------------------------------------------------------
# 6. Re-create the "_GoBack" bookmark right in the middle of the word
#    "pisane" (between "pis" and "ane"), matching where the cursor
#    was left after the edit.
# ------------------------------------------------------------------
$splitOffset = "Prerobit priority ciselne (strojove) na pis".Length
$bmPos = $newParaStart + $splitOffset
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
